# feat: Add customer loan count metrics + CAC per loan amortization
#
# Refresh of the cohort-performance export: fully_paid_count (and the
# metrics derived from it - fully_paid_rate, avg_ltv_to_cac_ratio, ...)
# reflect the updated loan-count/amortization logic, plus a handful of
# benign floating-point re-roundings that came along with the refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 177.4196705559369
$ws.Range("Q2").Value = 177.4196705559369
$ws.Range("T2").Value = 0.9985282363248826
$ws.Range("X2").Value = 1140
$ws.Range("AB2").Value = 0.7824296355247498
$ws.Range("AE2").Value = 1.203780779789366
$ws.Range("I3").Value = 7740.429999999998
$ws.Range("N3").Value = 7217.823
$ws.Range("X3").Value = 103
$ws.Range("AB3").Value = 1
$ws.Range("AE3").Value = 0.5430192140824556
$ws.Range("AF3").Value = 1.220323634979688
$ws.Range("E4").Value = 836425.4
$ws.Range("I4").Value = 84780.06999999999
$ws.Range("T4").Value = 0.9997304604276772
$ws.Range("X4").Value = 547
$ws.Range("AA4").Value = 716567.2300000001
$ws.Range("AB4").Value = 0.9579684734344482
$ws.Range("AE4").Value = 0.6544808264905551
$ws.Range("AF4").Value = 4.213231128770636
$ws.Range("X6").Value = 8
$ws.Range("AB6").Value = 1
$ws.Range("E7").Value = 2545370.78
$ws.Range("H7").Value = 496462.91
$ws.Range("I7").Value = 340922.6199999999
$ws.Range("O7").Value = 496462.91
$ws.Range("P7").Value = 495954.6266797292
$ws.Range("R7").Value = 209.4402984289397
$ws.Range("T7").Value = 0.9989761907485277
$ws.Range("X7").Value = 2003
$ws.Range("AB7").Value = 0.8458614945411682
$ws.Range("AD7").Value = 495954.6266797292
$ws.Range("AE7").Value = 1.454724483175093
$ws.Range("AF7").Value = 9.517958003190941
$ws.Range("K8").Value = 78.00227678571427
$ws.Range("Q8").Value = 78.00227678571427
$ws.Range("R8").Value = 78.00227678571427
$ws.Range("X8").Value = 224
$ws.Range("AB8").Value = 1
$ws.Range("AF8").Value = 1.939574450560139
$ws.Range("G9").Value = 1412.026268518518
$ws.Range("H9").Value = 231028.51
$ws.Range("N9").Value = 41511.60550000001
$ws.Range("O9").Value = 231028.51
$ws.Range("T9").Value = 0.9998682581341914
$ws.Range("X9").Value = 1062
$ws.Range("AB9").Value = 0.9833333492279053
$ws.Range("AE9").Value = 0.8382440498110526
$ws.Range("H10").Value = 813.3199999999999
$ws.Range("O10").Value = 813.3199999999999
$ws.Range("P10").Value = 813.3199999999999
$ws.Range("X10").Value = 4
$ws.Range("AB10").Value = 1
$ws.Range("AD10").Value = 813.3199999999999
$ws.Range("X11").Value = 18
$ws.Range("AB11").Value = 1
$ws.Range("X12").Value = 1936
$ws.Range("AB12").Value = 0.8938134908676147
$ws.Range("AE12").Value = 1.343816377676004
$ws.Range("N13").Value = 4709.8005
$ws.Range("X13").Value = 240
$ws.Range("AB13").Value = 1
$ws.Range("AF13").Value = 5.387279142715281
$ws.Range("I14").Value = 168572.54
$ws.Range("T14").Value = 0.99985397431846
$ws.Range("X14").Value = 1111
$ws.Range("AB14").Value = 0.9797177910804749
$ws.Range("AE14").Value = 0.553283102621172
$ws.Range("G15").Value = 906.2245454545454
$ws.Range("X15").Value = 22
$ws.Range("AB15").Value = 1
